$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# Update "想去人数" (want-to-go count) column F on 展览 sheet
$ws1.Range("F3").Value = 189
$ws1.Range("F5").Value = 5133
$ws1.Range("F9").Value = 570
$ws1.Range("F10").Value = 525
$ws1.Range("F13").Value = 1429
$ws1.Range("F14").Value = 4004
$ws1.Range("F15").Value = 424
$ws1.Range("F16").Value = 156
$ws1.Range("F17").Value = 140
$ws1.Range("F19").Value = 3082
$ws1.Range("F20").Value = 149
$ws1.Range("F21").Value = 1041
$ws1.Range("F25").Value = 82
$ws1.Range("F30").Value = 7
$ws1.Range("F33").Value = 4

# Update "想去人数" (want-to-go count) column F on 全部类型 sheet
$ws4.Range("F3").Value = 189
$ws4.Range("F6").Value = 5133
$ws4.Range("F10").Value = 570
$ws4.Range("F11").Value = 525
$ws4.Range("F14").Value = 1429
$ws4.Range("F15").Value = 4004
$ws4.Range("F16").Value = 424
$ws4.Range("F17").Value = 156
$ws4.Range("F18").Value = 140
$ws4.Range("F20").Value = 3082
$ws4.Range("F21").Value = 149
$ws4.Range("F22").Value = 1041
$ws4.Range("F26").Value = 82
$ws4.Range("F31").Value = 7
$ws4.Range("F34").Value = 4
